$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Insert a blank row before the old "USER FEATURES" header row (row 11),
# shifting rows 11-36 down to 12-37.
$ws.Rows("11").Insert()

# The content that used to live in the last row (old row 36, "Change
# Admin/Users so that it displays the user avatars") is moved up into the
# freshly-opened row 10, and gets a new "NEXT" status (bold, no fill).
$ws.Range("A10").Value = "Change Admin/Users so that it displays the user avatars "
$ws.Range("B10").Value = "Users, Global"
$ws.Range("C10").Value = "NEXT"
$ws.Range("C10").Font.Bold = $true
$ws.Range("D10").Value = "Do this in a way that it can be easily recycled on other pages. "

# Remove the now-duplicated old content (it shifted from row 36 to row 37
# when we inserted the row above).
$ws.Rows("37").Delete()

# Expand on the "UserGames model/table" note.
$ws.Range("D20").Value = "Requires some design. End goal is to support multiple games. Easier to do now rather than later.  Should give it a better than than UserGames. E.g. GameLogs"

# "Show User name, avatar image on index after login" is now Done.
$ws.Range("C23").Value = "Done"
$ws.Range("C23").Style = "Good"

# "User Partial View for insertion on multiple pages" now has a status.
$ws.Range("C24").Value = "Done"
$ws.Range("C24").Style = "Good"

# "Fill in content for legal page and hook up the link properly." is Done
# (legal screen hooked up).
$ws.Range("C36").Value = "Done"
$ws.Range("C36").Style = "Good"

# Update the saved selection to match the author's last cursor position.
[void]$ws.Activate()
[void]$ws.Range("C10").Select()
